# Sommerstudenter_timeliste.xlsx - "Add case 4 and 5"
# Fill in a full work week (MA-FR = 7.5 hours/day) for the two timesheet
# rows that previously had no hours logged (week rows at B17=28 and
# B23=31), then leave the selection where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (day "28"): Monday-Friday = 7.5 hours each -> SUM (M17) becomes 37.5
$ws.Range("F17:J17").Value = 7.5

# Row 23 (day "31"): Monday-Friday = 7.5 hours each -> SUM (M23) becomes 37.5
$ws.Range("F23:J23").Value = 7.5

# Scroll the sheet down a bit and leave the cursor on N15, matching
# where the author ended up when they saved.
$ws.Range("N15").Select()
